# casestudy.pptx — "Fixed comments and plots info"
$p = $ppt.ActivePresentation

# --- Slide 1: title text + subtitle position -------------------------------
$s1 = $p.Slides.Item(1)

# "po uzyskaniu dyplomu" -> "od uzyskania dyplomu" (2nd run of the title, after the line break)
$titleShape = $s1.Shapes.Item(1)
$titleTr = $titleShape.TextFrame.TextRange
$titleTr.Characters(57, 20).Text = "od uzyskania dyplomu"

# Move the subtitle placeholder up (y: 3750197 -> 2825963 EMU == 295.2911 -> 222.5168 pt)
$subtitleShape = $s1.Shapes.Item(2)
$subtitleShape.Top = 222.5168

# --- Slide 3: drop leading tab + remove the stray picture -------------------
$s3 = $p.Slides.Item(3)
$s3.Shapes.Item(2).TextFrame.TextRange.Characters(1, 1).Text = ""
$s3.Shapes.Item(4).Delete()

# --- Slide 4: drop leading tab ----------------------------------------------
$s4 = $p.Slides.Item(4)
$s4.Shapes.Item(2).TextFrame.TextRange.Characters(1, 1).Text = ""

# --- Slide 5: drop leading tab ----------------------------------------------
$s5 = $p.Slides.Item(5)
$s5.Shapes.Item(2).TextFrame.TextRange.Characters(1, 1).Text = ""

# --- Slide 6: drop leading tab + swap trailing comma for a period ----------
$s6 = $p.Slides.Item(6)
$s6.Shapes.Item(2).TextFrame.TextRange.Text = "Zarobki połowy absolwentów po pięciu latach nie przekraczają 4000 złotych. "
